# Add 2022-Q3 data.
#
# 1. Update the "总计" (totals) sheet: insert a new "2022-Q3" summary row
#    ahead of the existing "2022-Q2" / "2021-Q2" rows.
# 2. Insert a brand-new "2022-Q3" worksheet (positioned right after "总计",
#    before "2022-Q2") containing the per-fund holdings table for that
#    quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet - rewrite the 3-row data block as a 4-row block with the
#    new quarter inserted first.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totalsRows = @(
    @(0, "2022-Q3", 13, 0.8),
    @(1, "2022-Q2", 11, 1.17),
    @(2, "2021-Q2", 4, 0.74)
)

for ($i = 0; $i -lt $totalsRows.Length; $i++) {
    $r = $i + 2
    $row = $totalsRows[$i]

    $totals.Cells.Item($r, 1).Value = $row[0]
    $totals.Cells.Item($r, 2).Value = $row[1]
    $totals.Cells.Item($r, 3).Value = $row[2]
    $totals.Cells.Item($r, 4).Value = $row[3]
}

# Row 4 is brand new - give column A the same bold/centered/bordered style
# used by the rest of column A in this sheet.
$totals.Range("A2").Copy()
$totals.Range("A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet - clone the existing "2022-Q2" sheet (same
#    column layout/formatting), drop it in before "2022-Q2", rename it,
#    then replace its contents with the 2022-Q3 fund table.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundRows = @(
    @("010751", "宝盈优质成长混合A", "4.72", "92.28", "4.05", "0.1912", 8),
    @("012988", "嘉合锦明混合C", "3.37", "63.86", "3.53", "0.1190", 9),
    @("015011", "嘉合锦鑫混合C", "1.57", "62.86", "5.37", "0.0843", 3),
    @("006424", "嘉合锦程价值精选混合A", "1.47", "81.68", "5.38", "0.0791", 6),
    @("015010", "嘉合锦鑫混合A", "1.42", "62.86", "5.37", "0.0763", 3),
    @("012987", "嘉合锦明混合A", "2.15", "63.86", "3.53", "0.0759", 9),
    @("006425", "嘉合锦程价值精选混合C", "1.06", "81.68", "5.38", "0.0570", 6),
    @("233001", "大摩基础行业混合", "0.70", "78.50", "5.72", "0.0400", 8),
    @("010752", "宝盈优质成长混合C", "0.76", "92.28", "4.05", "0.0308", 8),
    @("005091", "嘉合睿金定期开放灵活配置混合型发起式C", "0.35", "74.56", "4.90", "0.0172", 7),
    @("011015", "嘉合锦元回报混合A", "0.78", "20.08", "1.85", "0.0144", 4),
    @("005090", "嘉合睿金定期开放灵活配置混合型发起式A", "0.26", "74.56", "4.90", "0.0127", 7),
    @("011016", "嘉合锦元回报混合C", "0.18", "20.08", "1.85", "0.0033", 4)
)

# The cloned sheet only has 11 data rows (rows 2-12); give column A the
# correct style for the two extra rows (13, 14) before filling values.
$q3.Range("A12").Copy()
$q3.Range("A13:A14").PasteSpecial(-4122)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q3.Cells.Item($r, 1).Value = $i

    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[0]

    $q3.Cells.Item($r, 3).NumberFormat = "@"
    $q3.Cells.Item($r, 3).Value = $row[1]

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[2]

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[3]

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[4]

    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[5]

    $q3.Cells.Item($r, 8).Value = $row[6]
}

# Restore the originally-active sheet (cloning/renaming moves focus to the
# freshly created sheet as a side effect).
$wb.Worksheets.Item("2021-Q2").Activate()
